$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: insert a new column at C, splitting "Strategy/Currency" into
# "Strategy" (col B) and "Override Currency" (new col C).
# ---------------------------------------------------------------------------
$ws1.Columns("C").Insert()

$ws1.Range("B1").Locked = $false
$ws1.Range("B1").Value = "Strategy"
$ws1.Range("D1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)

$ws1.Range("C1").Locked = $false
$ws1.Range("C1").Value = "Override Currency"
$ws1.Range("D1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet2: trim lookup column B down to 2 entries, and add a new lookup
# column C with currency codes.
# ---------------------------------------------------------------------------
$ws2.Range("B3").ClearContents()
$ws2.Range("B4").ClearContents()
$ws2.Range("B5").ClearContents()

$ws2.Range("C1").Value = "CAD"
$ws2.Range("C2").Value = "EUR"
$ws2.Range("C3").Value = "GBP"
$ws2.Range("C4").Value = "JPY"
$ws2.Range("C5").Value = "USD"
$ws2.Range("B1").Copy()
$ws2.Range("C1:C5").PasteSpecial(-4122)
$ws2.Range("C1").Value = "CAD"
$ws2.Range("C2").Value = "EUR"
$ws2.Range("C3").Value = "GBP"
$ws2.Range("C4").Value = "JPY"
$ws2.Range("C5").Value = "USD"

$wb.Save()
